$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hedge")

# Insert a new row at row 66 - existing row 66 and below shift down to 67+
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new fund entry (Braidwell / Alex Carnal)
$ws.Range("B66").Value = "Braidwell"
$ws.Range("C66").Value = "Alex Carnal"
$ws.Range("P66").Value = 2722.6707569999999
$ws.Range("Q66").Value = 3102.7736410000002
$ws.Range("R66").Value = 3187.7679090000001
$ws.Range("S66").Value = "Healthcare"

# Hyperlinks matching the SEC EDGAR filer pattern used throughout the sheet
$ws.Hyperlinks.Add($ws.Range("B66"), "https://www.sec.gov/edgar/browse/?CIK=1920938") | Out-Null
$ws.Hyperlinks.Add($ws.Range("P66"), "https://www.sec.gov/Archives/edgar/data/1920938/000142050624000478/xslForm13F_X02/primary_doc.xml") | Out-Null
$ws.Hyperlinks.Add($ws.Range("Q66"), "https://www.sec.gov/Archives/edgar/data/1920938/000192093824000004/xslForm13F_X02/primary_doc.xml") | Out-Null
$ws.Hyperlinks.Add($ws.Range("R66"), "https://www.sec.gov/Archives/edgar/data/1920938/000142050624001566/xslForm13F_X02/primary_doc.xml") | Out-Null

# Leave cursor/selection on the newly entered cell, matching the author's editing session
$ws.Range("S66").Select()
